# Update dSF column (F) values with newly repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 0
